# Update the "Last Updated" timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 01:21 PM"

# The "Stock List" sheet refreshed: the former first data row (CAPTRU-RE1)
# dropped off the list, every other row shifted up by one, and a new row
# (TRAVELFOOD) was appended at the bottom.
$ws = $wb.Worksheets.Item("Stock List")
$ws.Rows.Item(2).Delete()

$ws.Range("A76").Value = "📋"
$ws.Range("B76").Value = "TRAVELFOOD"
$ws.Range("C76").Value = "TRAVELFOOD"
$ws.Range("D76").Value = 1316.3
$ws.Range("E76").Value = 0.1141
$ws.Range("F76").Value = "N/A"
$ws.Range("G76").Value = "N/A"
$ws.Range("H76").Value = 17332.9705
